$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('A7').Value2 = 'Prabowo Instruksikan Perbaikan MBG: Minta Koki Terlatih hingga Pasang CCTV'
$ws.Range('B7').Value2 = '2025-09-29T09:06:44+07:00'
$ws.Range('C7').Value2 = 'Anggi Muliawati'
$ws.Range('D7').Value2 = 'https://www.detik.com/bali/berita/d-8135103/prabowo-instruksikan-perbaikan-mbg-minta-koki-terlatih-hingga-pasang-cctv'

$ws.Range('A8').Value2 = 'Sederet Instruksi Prabowo ke BGN Buntut Marak Keracunan MBG'
$ws.Range('B8').Value2 = '2025-09-29T08:43:05+07:00'
$ws.Range('C8').Value2 = 'Anggi Muliawati'
$ws.Range('D8').Value2 = 'https://www.detik.com/jateng/berita/d-8135088/sederet-instruksi-prabowo-ke-bgn-buntut-marak-keracunan-mbg'

$ws.Range('A9').Value2 = 'Legislator Setuju Koki Makan Bergizi Gratis Harus Punya Pengalaman'
$ws.Range('B9').Value2 = '2025-09-29T08:38:02+07:00'
$ws.Range('C9').Value2 = 'Isal Mawardi'
$ws.Range('D9').Value2 = 'https://news.detik.com/berita/d-8135085/legislator-setuju-koki-makan-bergizi-gratis-harus-punya-pengalaman'

$ws.Range('A10').Value2 = 'Juru Masak MBG Dievaluasi'
$ws.Range('B10').Value2 = '2025-09-29T07:58:00+07:00'
$ws.Range('C10').Value2 = 'Trypama Randra'
$ws.Range('D10').Value2 = 'https://news.detik.com/berita/d-8134900/juru-masak-mbg-dievaluasi'

$ws.Range('A11').Value2 = 'Terungkap Bakteri ''Biang Kerok'' Keracunan MBG di Bandung Barat'
$ws.Range('B11').Value2 = '2025-09-29T07:01:12+07:00'
$ws.Range('C11').Value2 = 'Antara'
$ws.Range('D11').Value2 = 'https://news.detik.com/berita/d-8135037/terungkap-bakteri-biang-kerok-keracunan-mbg-di-bandung-barat'

$ws.Range('A12').Value2 = 'Bingkai Sepekan: MBG Jadi Polemik, Menu Minim Gizi-Ribuan Keracunan'
$ws.Range('B12').Value2 = '2025-09-27T07:00:56+07:00'
$ws.Range('C12').Value2 = '-'
$ws.Range('D12').Value2 = 'https://health.detik.com/fotohealth/d-8132256/bingkai-sepekan-mbg-jadi-polemik-menu-minim-gizi-ribuan-keracunan'

$ws.Range('A13').Value2 = '''''Tur'' ke Dapur Makan Bergizi Gratis di Bandung'
$ws.Range('B13').Value2 = '2025-09-26T19:00:53+07:00'
$ws.Range('C13').Value2 = '-'
$ws.Range('D13').Value2 = 'https://www.detik.com/jabar/foto/d-8131846/tur-ke-dapur-makan-bergizi-gratis-di-bandung'

$ws.Range('A14').Value2 = 'MBG di Sekolah Jakut Isinya Cuma Snack, Gizi Dipertanyakan'
$ws.Range('B14').Value2 = '2025-09-26T17:00:17+07:00'
$ws.Range('C14').Value2 = '-'
$ws.Range('D14').Value2 = 'https://health.detik.com/fotohealth/d-8131726/mbg-di-sekolah-jakut-isinya-cuma-snack-gizi-dipertanyakan'

$ws.Range('A15').Value2 = 'Penampakan Dapur MBG di Bandung yang Bikin Ratusan Siswa Keracunan'
$ws.Range('B15').Value2 = '2025-09-25T18:00:53+07:00'
$ws.Range('C15').Value2 = '-'
$ws.Range('D15').Value2 = 'https://news.detik.com/foto-news/d-8130038/penampakan-dapur-mbg-di-bandung-yang-bikin-ratusan-siswa-keracunan'

$ws.Range('A16').Value2 = 'Potret Pilu Ribuan Anak Sekolah Jadi Korban Makan Bergizi Gratis'
$ws.Range('B16').Value2 = '2025-09-25T13:10:51+07:00'
$ws.Range('C16').Value2 = '-'
$ws.Range('D16').Value2 = 'https://health.detik.com/fotohealth/d-8129561/potret-pilu-ribuan-anak-sekolah-jadi-korban-makan-bergizi-gratis'

$ws.Range('A17').Value2 = '4 Instruksi Prabowo Saat Panggil Kepala BGN Usai Marak Keracunan MBG'
$ws.Range('B17').Value2 = '2025-09-29T06:29:26+07:00'
$ws.Range('C17').Value2 = 'Anggi Muliawati'
$ws.Range('D17').Value2 = 'https://news.detik.com/berita/d-8135031/4-instruksi-prabowo-saat-panggil-kepala-bgn-usai-marak-keracunan-mbg'

$ws.Range('A18').Value2 = 'CISDI Dorong Pemerintah Susun Perpres Perbaikan Menyeluruh MBG'
$ws.Range('B18').Value2 = '2025-09-29T06:20:00+07:00'
$ws.Range('C18').Value2 = 'Isal Mawardi'
$ws.Range('D18').Value2 = 'https://news.detik.com/berita/d-8135020/cisdi-dorong-pemerintah-susun-perpres-perbaikan-menyeluruh-mbg'

$ws.Range('A19').Value2 = 'Langkah Prabowo untuk Kelanjutan Program MBG'
$ws.Range('B19').Value2 = '2025-09-29T05:58:59+07:00'
$ws.Range('C19').Value2 = 'Tim detikcom'
$ws.Range('D19').Value2 = 'https://www.detik.com/kalimantan/berita/d-8134967/langkah-prabowo-untuk-kelanjutan-program-mbg'

$ws.Range('A20').Value2 = 'Blak-blakan Menkop soal Tantangan Koperasi Merah Putih'
$ws.Range('B20').Value2 = '2025-09-26T17:03:11+07:00'
$ws.Range('C20').Value2 = 'Wisma Putra'
$ws.Range('D20').Value2 = 'https://www.detik.com/jabar/bisnis/d-8132102/blak-blakan-menkop-soal-tantangan-koperasi-merah-putih'

$ws.Range('A21').Value2 = 'Begini Strategi Khofifah Cegah Jeratan Pinjol Ilegal di Tengah Warga'
$ws.Range('B21').Value2 = '2025-09-26T10:45:16+07:00'
$ws.Range('C21').Value2 = 'Faiq Azmi'
$ws.Range('D21').Value2 = 'https://www.detik.com/jatim/berita/d-8131171/begini-strategi-khofifah-cegah-jeratan-pinjol-ilegal-di-tengah-warga'

$ws.Range('A22').Value2 = 'Kemenkop Puji Kopdes Aeng Batu-batu Takalar Raup Omzet Rp 2 Juta Per Hari'
$ws.Range('B22').Value2 = '2025-09-23T14:46:31+07:00'
$ws.Range('C22').Value2 = 'Adhe Junaedi Sholat'
$ws.Range('D22').Value2 = 'https://www.detik.com/sulsel/bisnis/d-8126210/kemenkop-puji-kopdes-aeng-batu-batu-takalar-raup-omzet-rp-2-juta-per-hari'

$ws.Range('A23').Value2 = 'Kemenkop Kawal Pemberian Pinjaman Modal 3.059 Koperasi Merah Putih di Sulsel'
$ws.Range('B23').Value2 = '2025-09-23T11:26:57+07:00'
$ws.Range('C23').Value2 = 'Adhe Junaedi Sholat'
$ws.Range('D23').Value2 = 'https://www.detik.com/sulsel/bisnis/d-8125752/kemenkop-kawal-pemberian-pinjaman-modal-3-059-koperasi-merah-putih-di-sulsel'

$ws.Range('A24').Value2 = 'Ujian Awal Koperasi Merah Putih di Usia Setipis Daun Kelor'
$ws.Range('B24').Value2 = '2025-09-22T10:30:00+07:00'
$ws.Range('C24').Value2 = 'Andry Haryanto'
$ws.Range('D24').Value2 = 'https://www.detik.com/jabar/berita/d-8123636/ujian-awal-koperasi-merah-putih-di-usia-setipis-daun-kelor'

$ws.Range('A25').Value2 = 'PPPK Bisa Diperbantukan di Kopdes, Pinjaman Rp 3 M Cair Pekan Depan'
$ws.Range('B25').Value2 = '2025-09-18T19:25:18+07:00'
$ws.Range('C25').Value2 = 'Agus Setyadi'
$ws.Range('D25').Value2 = 'https://www.detik.com/sumut/bisnis/d-8119042/pppk-bisa-diperbantukan-di-kopdes-pinjaman-rp-3-m-cair-pekan-depan'

$ws.Range('A26').Value2 = 'Melihat Koperasi Merah Putih di Melawai Jaksel'
$ws.Range('B26').Value2 = '2025-07-23T13:00:24+07:00'
$ws.Range('C26').Value2 = '-'
$ws.Range('D26').Value2 = 'https://finance.detik.com/foto-bisnis/d-8024832/melihat-koperasi-merah-putih-di-melawai-jaksel'

$ws.Range('A27').Value2 = 'Wamenkop Apresiasi 100% Pembentukan Badan Hukum Kopdes Merah Putih di Lahat'
$ws.Range('B27').Value2 = '2025-06-10T15:10:00+07:00'
$ws.Range('C27').Value2 = '-'
$ws.Range('D27').Value2 = 'https://www.detik.com/sumbagsel/foto/d-7957338/wamenkop-apresiasi-100-pembentukan-badan-hukum-kopdes-merah-putih-di-lahat'

$ws.Range('A28').Value2 = 'Momen Wamenkop Ferry Pantau Musdesus Kopdes Merah Putih di Padang'
$ws.Range('B28').Value2 = '2025-05-29T19:30:25+07:00'
$ws.Range('C28').Value2 = '-'
$ws.Range('D28').Value2 = 'https://www.detik.com/sumut/foto/d-7939188/momen-wamenkop-ferry-pantau-musdesus-kopdes-merah-putih-di-padang'

$ws.Range('A29').Value2 = 'Budi Arie Sambangi KPK, Minta Pengawalan Program Koperasi Desa'
$ws.Range('B29').Value2 = '2025-05-21T17:29:29+07:00'
$ws.Range('C29').Value2 = '-'
$ws.Range('D29').Value2 = 'https://news.detik.com/foto-news/d-7925409/budi-arie-sambangi-kpk-minta-pengawalan-program-koperasi-desa'

$ws.Range('A30').Value2 = '20 Ribu Kopdes Ditargetkan Dapat Pinjaman dari Bank BUMN Tahun Ini'
$ws.Range('B30').Value2 = '2025-09-18T13:53:03+07:00'
$ws.Range('C30').Value2 = 'Retno Ayuningrum'
$ws.Range('D30').Value2 = 'https://finance.detik.com/moneter/d-8118291/20-ribu-kopdes-ditargetkan-dapat-pinjaman-dari-bank-bumn-tahun-ini'

$ws.Range('A31').Value2 = 'Link Pengumuman Hasil Seleksi Administrasi PMO Koperasi Merah Putih Kemenkop'
$ws.Range('B31').Value2 = '2025-09-15T14:30:00+07:00'
$ws.Range('C31').Value2 = 'Nikita Rosa'
$ws.Range('D31').Value2 = 'https://www.detik.com/edu/detikpedia/d-8112624/link-pengumuman-hasil-seleksi-administrasi-pmo-koperasi-merah-putih-kemenkop'

$ws.Range('A32').Value2 = 'Zulhas Minta Percepatan Pinjaman Kopdes Merah Putih: Rp 200 T Ada di Bank Himbara'
$ws.Range('B32').Value2 = '2025-09-15T12:18:25+07:00'
$ws.Range('C32').Value2 = 'Lisye Sri Rahayu'
$ws.Range('D32').Value2 = 'https://news.detik.com/berita/d-8112520/zulhas-minta-percepatan-pinjaman-kopdes-merah-putih-rp-200-t-ada-di-bank-himbara'

$ws.Range('A33').Value2 = 'Menteri Koperasi Tinjau Koperasi Merah Putih di Tuban, Ini Hasilnya'
$ws.Range('B33').Value2 = '2025-09-13T19:45:58+07:00'
$ws.Range('C33').Value2 = 'Ainur Rofiq'
$ws.Range('D33').Value2 = 'https://www.detik.com/jatim/bisnis/d-8110591/menteri-koperasi-tinjau-koperasi-merah-putih-di-tuban-ini-hasilnya'
